$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "318.04"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "3.62%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "39.97"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.82%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.138"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.66%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08231"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.82%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.084"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "7.34%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.336"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "4.69%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9417"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.17%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1367"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-6.80%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1980"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.70%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09102"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.84%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03521"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.57%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09811"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.19%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001386"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.15%"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.04346"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.53%"
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006325"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "8.32%"
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.702"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.25%"
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.322"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "3.34%"
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.242"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-6.25%"
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3500"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.14%"
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1309"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.53%"
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "MCDex"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.986"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "9.28%"
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.2443"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.24%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001226"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.81%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004835"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "12.86%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001296"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.37%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02212"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "8.31%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05229"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.00%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007681"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.09%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009727"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-5.54%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "4.42%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002035"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-4.07%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008903"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-1.97%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006604"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "6.41%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.37%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002929"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-5.67%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.37%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.37%"
